$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace color labels in column B: "blue" -> "deepskyblue", "black" -> "gold"
$rng = $ws.Range("B1:B16")
$rng.Replace("blue", "deepskyblue", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$rng.Replace("black", "gold", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Update the active selection to match the authored state
$ws.Range("D8").Select()
